# Scheduled runner update: refresh derived profit/price columns (H, I-N)
# across the Tonberry_Profits leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below are the latest market-board snapshot pulled by the scraper;
# a handful of rows also gain/lose a trailing N (or M) column cell depending on
# whether the HQ/NQ leve-turn-in math now yields a value for that row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1063.3334
$ws.Range("J17").Value = 893.8475
$ws.Range("L17").Value = 2681.5425
$ws.Range("N17").Value = -3017.5425
$ws.Range("H40").Value = 2347.1428
$ws.Range("I40").Value = 2420
$ws.Range("J40").Value = 2250
$ws.Range("K40").Value = 2420
$ws.Range("L40").Value = 2250
$ws.Range("M40").Value = -2245
$ws.Range("N40").Value = -2600
$ws.Range("H69").Value = 1509.8
$ws.Range("I69").Value = 1509.8
$ws.Range("K69").Value = 4529.4
$ws.Range("M69").Value = -3655.4
$ws.Range("H72").Value = 1509.8
$ws.Range("I72").Value = 1509.8
$ws.Range("K72").Value = 13588.2
$ws.Range("M72").Value = -9220.199999999999
$ws.Range("H88").Value = 2411.5
$ws.Range("J88").Value = 1383.3334
$ws.Range("L88").Value = 1383.3334
$ws.Range("N88").Value = -2195.3334
$ws.Range("H91").Value = 2411.5
$ws.Range("J91").Value = 1383.3334
$ws.Range("L91").Value = 1383.3334
$ws.Range("N91").Value = -4191.3334
$ws.Range("H100").Value = 1531.8
$ws.Range("I100").Value = 1368.6666
$ws.Range("K100").Value = 1368.6666
$ws.Range("M100").Value = -827.6666
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()  # was -84987
$ws.Range("H125").Value = 690.8889
$ws.Range("I125").Value = 369.66666
$ws.Range("K125").Value = 3326.99994
$ws.Range("M125").Value = -866.9999399999997
$ws.Range("H132").Value = 1105.5555
$ws.Range("J132").Value = 1624.25
$ws.Range("L132").Value = 4872.75
$ws.Range("N132").Value = -9932.75
$ws.Range("H138").Value = 2823.963
$ws.Range("I138").Value = 2844
$ws.Range("K138").Value = 8532
$ws.Range("M138").Value = -3392
$ws.Range("H141").Value = 1477127.8
$ws.Range("I141").Value = 2156257
$ws.Range("K141").Value = 6468771
$ws.Range("M141").Value = -6463591

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1728.2222
$ws.Range("I45").Value = 1602
$ws.Range("J45").Value = 1854.4445
$ws.Range("K45").Value = 1602
$ws.Range("L45").Value = 1854.4445
$ws.Range("M45").Value = -1225
$ws.Range("N45").Value = -2608.4445
$ws.Range("H88").Value = 4699.7
$ws.Range("J88").Value = 5312.25
$ws.Range("L88").Value = 5312.25
$ws.Range("N88").Value = -6124.25
$ws.Range("H91").Value = 4699.7
$ws.Range("J91").Value = 5312.25
$ws.Range("L91").Value = 5312.25
$ws.Range("N91").Value = -8120.25
$ws.Range("H122").Value = 1506
$ws.Range("I122").Value = 1506
$ws.Range("K122").Value = 4518
$ws.Range("M122").Value = -2068

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 824.5
$ws.Range("I22").Value = 824.5
$ws.Range("K22").Value = 824.5
$ws.Range("M22").Value = -651.5
$ws.Range("H94").Value = 893.4
$ws.Range("I94").Value = 859.3333
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 859.3333
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -408.3333
$ws.Range("N94").Value = -2102
$ws.Range("H137").Value = 62000
$ws.Range("J137").Value = 62000
$ws.Range("L137").Value = 62000
$ws.Range("N137").Value = -72200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()  # was -38195
$ws.Range("H132").Value = 2002.9642
$ws.Range("I132").Value = 1140.6842
$ws.Range("J132").Value = 3823.3333
$ws.Range("K132").Value = 3422.0526
$ws.Range("L132").Value = 11469.9999
$ws.Range("M132").Value = -892.0526
$ws.Range("N132").Value = -16529.9999
$ws.Range("H134").Value = 747.6539
$ws.Range("I134").Value = 747.6539
$ws.Range("K134").Value = 2242.9617
$ws.Range("M134").Value = 292.0383000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 125001830
$ws.Range("J116").Value = 166668670
$ws.Range("L116").Value = 500006010
$ws.Range("N116").Value = -500012894
$ws.Range("H125").Value = 1430
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()  # was -54840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1837.5
$ws.Range("I122").Value = 1527.25
$ws.Range("K122").Value = 4581.75
$ws.Range("M122").Value = -2131.75
$ws.Range("H123").Value = 15682.429
$ws.Range("J123").Value = 15682.429
$ws.Range("L123").Value = 15682.429
$ws.Range("N123").Value = -20582.429
$ws.Range("H132").Value = 1751872.6
$ws.Range("I132").Value = 2139547.5
$ws.Range("K132").Value = 6418642.5
$ws.Range("M132").Value = -6416112.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6660.909
$ws.Range("J7").Value = 8280.857
$ws.Range("L7").Value = 8280.857
$ws.Range("N7").Value = -8504.857
$ws.Range("H40").Value = 8489.6
$ws.Range("I40").Value = 1979.8
$ws.Range("J40").Value = 14999.4
$ws.Range("K40").Value = 1979.8
$ws.Range("L40").Value = 14999.4
$ws.Range("M40").Value = -1843.8
$ws.Range("N40").Value = -15271.4
$ws.Range("H122").Value = 13582.333
$ws.Range("I122").Value = 12873.5
$ws.Range("K122").Value = 38620.5
$ws.Range("M122").Value = -36170.5
$ws.Range("H126").Value = 6660.909
$ws.Range("J126").Value = 8280.857
$ws.Range("L126").Value = 24842.571
$ws.Range("N126").Value = -29782.571
$ws.Range("H132").Value = 3424.5
$ws.Range("I132").Value = 2566.3333
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 7698.999899999999
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -5168.999899999999
$ws.Range("N132").Value = -23057
$ws.Range("H134").Value = 49154
$ws.Range("J134").Value = 49154
$ws.Range("L134").Value = 49154
$ws.Range("N134").Value = -59294
$ws.Range("H136").Value = 4164.467
$ws.Range("I136").Value = 1867.5
$ws.Range("J136").Value = 4999.727
$ws.Range("K136").Value = 5602.5
$ws.Range("L136").Value = 14999.181
$ws.Range("M136").Value = -3052.5
$ws.Range("N136").Value = -20099.181
$ws.Range("H139").Value = 45000
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()  # was -34860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 46685.332
$ws.Range("J40").Value = 60028
$ws.Range("L40").Value = 60028
$ws.Range("N40").Value = -60326
$ws.Range("H95").Value = 47342.332
$ws.Range("J95").Value = 47342.332
$ws.Range("L95").Value = 47342.332
$ws.Range("N95").Value = -52834.332
$ws.Range("H132").Value = 3536.4614
$ws.Range("I132").Value = 2998.5715
$ws.Range("J132").Value = 4164
$ws.Range("K132").Value = 8995.7145
$ws.Range("L132").Value = 12492
$ws.Range("M132").Value = -6465.7145
$ws.Range("N132").Value = -17552
